$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
}

Replace-Text "2024-07-01 Monday" "2024-07-02 Tuesday"

Replace-Text "790÷8=98, 6" "337÷2=168, 1"
Replace-Text "142÷2=71, 0" "268÷7=38, 2"
Replace-Text "180÷9=20, 0" "810÷4=202, 2"
Replace-Text "132÷2=66, 0" "746÷2=373, 0"
Replace-Text "332÷3=110, 2" "129÷9=14, 3"

Replace-Text "829÷6=138, 1" "346÷2=173, 0"
Replace-Text "378÷2=189, 0" "331÷7=47, 2"
Replace-Text "516÷3=172, 0" "232÷9=25, 7"
Replace-Text "275÷8=34, 3" "956÷8=119, 4"
Replace-Text "981÷7=140, 1" "404÷2=202, 0"

Replace-Text "225÷7=32, 1" "922÷6=153, 4"
Replace-Text "452÷6=75, 2" "768÷4=192, 0"
Replace-Text "886÷9=98, 4" "302÷8=37, 6"
Replace-Text "845÷7=120, 5" "782÷6=130, 2"
Replace-Text "167÷9=18, 5" "939÷8=117, 3"

Replace-Text "786÷6=131, 0" "261÷8=32, 5"
Replace-Text "642÷2=321, 0" "209÷5=41, 4"
Replace-Text "761÷5=152, 1" "318÷5=63, 3"
Replace-Text "613÷2=306, 1" "844÷3=281, 1"
Replace-Text "963÷2=481, 1" "881÷9=97, 8"

Replace-Text "660÷5=132, 0" "444÷9=49, 3"
Replace-Text "718÷8=89, 6" "699÷6=116, 3"
Replace-Text "687÷4=171, 3" "744÷7=106, 2"
Replace-Text "584÷5=116, 4" "103÷2=51, 1"
Replace-Text "915÷9=101, 6" "176÷4=44, 0"
